# Inclusão de resultados de importância das variáveis
#
# Corrige os nomes das variáveis tratadas (coluna "Nome Variável Tratada")
# para que o prefixo (CAT_/NUM_/BIN_) seja coerente com o "Tipo Tratada"
# de cada linha, e atualiza a seleção ativa da planilha.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Q024 "Presença de computador" é Numérica -> deve ser NUM_Q024 (estava BIN_Q024)
$ws.Range("C45").Value = "NUM_Q024"
# Q025 "Acesso à internet" é Binária -> deve ser BIN_Q025 (estava NUM_Q025)
$ws.Range("C46").Value = "BIN_Q025"

# Remove o prefixo redundante "TP_" do nome tratado da dependência administrativa
$ws.Range("C11").Value = "CAT_DEPENDENCIA_ADM_ESC"
# Padroniza o sufixo do nome tratado da localização da escola
$ws.Range("C12").Value = "CAT_LOCALIZACAO_ESC"

# Atualiza a célula selecionada ativa na planilha
$ws.Range("F25").Select() | Out-Null
